$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Layout tweak: make room for a new "Ngành học" column in the header row (row 6) ---
# Before : A=MSSV  B=Tên SV  C=Giới tính  D=Ngày sinh  E=Khung chương trình   F=Học kỳ
# After  : A=MSSV  B=Tên SV  C=Ngày sinh  D=Giới tính  E=Ngành học (NEW)     F=Khung chương trình   G=Học kỳ

# Extend the header formatting (border/bold/centered style) from F6 into the new G6 cell
$ws.Range("F6").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Học kỳ moves one column to the right, Khung chương trình shifts into F,
# a brand-new Ngành học column is inserted at E, and Giới tính / Ngày sinh swap places
$ws.Range("G6").Value2 = "Học kỳ"
$ws.Range("F6").Value2 = "Khung chương trình"
$ws.Range("E6").Value2 = "Ngành học"
$ws.Range("D6").Value2 = "Giới tính"
$ws.Range("C6").Value2 = "Ngày sinh"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20.8776041667   # A: 31            -> 21.7109375
$ws.Columns.Item(2).ColumnWidth = 29.0221354167   # B: 14.5703125    -> 29.85546875
$ws.Columns.Item(4).ColumnWidth = 25.5924479167   # D: 30.42578125   -> 26.42578125
$ws.Columns.Item(7).ColumnWidth = 21.5924479167   # G: (new column)  -> 22.42578125

# --- Leave the selection where the editor left it ---
$ws.Range("E6").Select() | Out-Null
